$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: did not find -> $old"
    }
}

Replace-Text "2024-08-03 Saturday" "2024-08-04 Sunday"
Replace-Text "94×24=" "33×55="
Replace-Text "98×92=" "77×68="
Replace-Text "91×86=" "46×28="
Replace-Text "13×63=" "63×51="
Replace-Text "30×55=" "90×69="
Replace-Text "38×46=" "86×28="
Replace-Text "87×28=" "22×38="
Replace-Text "86×59=" "28×87="
Replace-Text "51×33=" "87×36="
Replace-Text "99×91=" "67×62="
Replace-Text "52×73=" "56×88="
Replace-Text "27×45=" "79×52="
Replace-Text "91×91=" "98×69="
Replace-Text "49×41=" "15×91="
Replace-Text "58×64=" "54×82="
Replace-Text "96×99=" "55×13="
Replace-Text "51×77=" "60×58="
Replace-Text "43×72=" "28×98="
Replace-Text "44×86=" "85×85="
Replace-Text "38×79=" "26×49="
Replace-Text "23×60=" "76×78="
Replace-Text "82×74=" "22×75="
Replace-Text "39×60=" "13×12="
Replace-Text "76×23=" "33×23="
Replace-Text "85×70=" "80×62="

Write-Host "Done."
